$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: add Name (column C) "Close to sheets by 840" and update Lat/Long
$ws.Range("C9").Value = "Close to sheets by 840"
$ws.Range("E9").Value = 36.0910691
$ws.Range("F9").Value = -79.7061625

# Row 10: remove Name (column C) - it was moved to row 9
$ws.Range("C10").ClearContents()

# Row 14: update Lat/Long
$ws.Range("E14").Value = 36.0599295
$ws.Range("F14").Value = -79.7151021

# Update selection to match final active cell F14
$ws.Range("F14").Select()
